$d = $word.ActiveDocument

# ------------------------------------------------------------------
# This document had:
#   ...paragraph ending "...Once they'"
#   an empty paragraph (firstLine-indent only, no runs)
#   paragraph starting "The goal for my senior project..." and ending
#       "...Users will have a map " + the hidden _GoBack bookmark
#
# The edit:
#   1. Move the (hidden) _GoBack bookmark from the end of the
#      "Users will have a map " paragraph to the end of the
#      "...Once they'" paragraph (collapsed, right after the text).
#   2. Delete the empty paragraph that sits between the two
#      paragraphs above (its mark only - no visible text is lost).
# ------------------------------------------------------------------

# Locate the paragraph that ends with "...Once they'" by scanning for
# its distinctive trailing text instead of hard-coding an index.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "*Once they*") {
        $targetIndex = $i
        break
    }
}

$p = $d.Paragraphs($targetIndex)

# Paragraph.Range.Text always carries a trailing paragraph-mark (\r),
# so the real end-of-content boundary is End-1, not End.
$contentEnd = $p.Range.End - 1

# Placing a *collapsed* Range exactly on a paragraph boundary confuses
# Bookmarks.Add in this host, so we temporarily insert a private-use
# placeholder character right at the boundary, anchor the bookmark
# just in front of it (now a safe, non-boundary position), and then
# remove the placeholder again.
$placeholderRange = $d.Range($contentEnd, $contentEnd)
$placeholderRange.InsertAfter([char]0xE000)

$bookmarkRange = $d.Range($contentEnd, $contentEnd)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

$cleanupRange = $d.Range($contentEnd, $contentEnd + 1)
$cleanupRange.Delete()

# Now remove the empty paragraph that immediately follows; re-fetch it
# fresh since the paragraph collection / offsets above may have shifted.
$emptyPara = $d.Paragraphs($targetIndex + 1)
$emptyRange = $d.Range($emptyPara.Range.Start, $emptyPara.Range.End)
$emptyRange.Delete()
